$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 06:44"

# Honduras overtook Austria in total cases, so the two countries swap
# positions between row 56 and row 57 (labels swap, values updated).
$ws.Range("A56").Value = "Honduras"
$ws.Range("A57").Value = "Austria"

# Fiyi and Dominica swap positions between row 205 and row 206
# (their underlying numbers are identical, only the labels swap).
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

# Row 15 - Pakistan: refreshed totals
$ws.Range("B15").Value = 206512
$ws.Range("C15").Value = 3557
$ws.Range("D15").Value = 95407
$ws.Range("E15").Value = 106938
$ws.Range("G15").Value = 49
$ws.Range("H15").Value = 4167

# Row 54 - Kazajistan: refreshed totals
$ws.Range("B54").Value = 21327
$ws.Range("C54").Value = 547
$ws.Range("D54").Value = 12933
$ws.Range("E54").Value = 8216

# Row 56 - now Honduras: refreshed totals
$ws.Range("B56").Value = 18082
$ws.Range("C56").Value = 1075
$ws.Range("D56").Value = 1875
$ws.Range("E56").Value = 15728
$ws.Range("H56").Value = 479

# Row 57 - now Austria: refreshed totals
$ws.Range("B57").Value = 17654
$ws.Range("D57").Value = 16401
$ws.Range("E57").Value = 551
$ws.Range("H57").Value = 702

# Row 83 - Haiti: refreshed totals
$ws.Range("B83").Value = 5847
$ws.Range("C83").Value = 70
$ws.Range("D83").Value = 787
$ws.Range("E83").Value = 4956
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 104

# Row 88 - Kirguistan: refreshed totals
$ws.Range("B88").Value = 5017
$ws.Range("C88").Value = 269
$ws.Range("D88").Value = 2294
$ws.Range("E88").Value = 2673
$ws.Range("G88").Value = 3
$ws.Range("H88").Value = 50

# Row 167 - Mongolia: refreshed totals
$ws.Range("B167").Value = 220
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 45

# Row 186 - Butan: refreshed totals
$ws.Range("D186").Value = 44
$ws.Range("E186").Value = 32
